$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "08 March 2025" -> "11" + " March 2025" (two runs) in the Date
# value cell of the first table.
# ---------------------------------------------------------------------------
$table1 = $d.Tables.Item(1)
$dateCell = $table1.Cell(1, 2)
$dateXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00C5158A" w:rsidRDefault="005109C2"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>11</w:t></w:r><w:r><w:t xml:space="preserve"> March 2025</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$dateCell.Range.InsertXML($dateXml)

# ---------------------------------------------------------------------------
# Change 2: Add a "4 Marks" run (Times New Roman, 12pt) before the _GoBack
# bookmark in the Maximum Marks value cell of the first table.
# ---------------------------------------------------------------------------
$marksCell = $table1.Cell(4, 2)
$marksXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00C5158A" w:rsidRDefault="00C5158A"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>4 Marks</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$marksCell.Range.InsertXML($marksXml)

# ---------------------------------------------------------------------------
# Change 3: Merge "March 1, 2025 - March 8, 202" + "5" into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("March 1, 2025 - March 8, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "March 1, 2025 - March 8, 2025", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: Merge "Streaming Quality & " + "Performance Testing" into a
# single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Streaming Quality & Performance Testing", $true, $false, $false, $false, $false, $true, 1, $false, "Streaming Quality & Performance Testing", 2) | Out-Null
